$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peer  and self assessment")

$ws.Range("B6").Value = "Insufficient"
$ws.Range("C6").Value = "Yet to participate in group meetings, yet to respond to ANY communication from anyone. Lukasz has also yet to respond. Quite frankly disgraceful unless a proper explanation as to why is given."

$ws.Range("B19").Value = "Insufficient"
$ws.Range("C19").Value = "UTP works on a person by person basis instead of in workshop groups. This might be one of the root causes of her lack of involvement inside the project. She has clearly shown a lack of motivation to work inside the group by refusing any and all communication. "
